$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 51.08364133333333
$ws.Cells.Item(2, 8).Value = 153.250924
$ws.Cells.Item(2, 9).Value = 0.2324760785757066
$ws.Cells.Item(2, 10).Value = 0.2324760785757066
$ws.Cells.Item(2, 13).Value = 17.367416
$ws.Cells.Item(2, 14).Value = 52.102248
$ws.Cells.Item(2, 15).Value = 0.199925246532591
$ws.Cells.Item(2, 16).Value = 0.199925246532591
$ws.Cells.Item(2, 17).Value = 887.1908498307946
$ws.Cells.Item(2, 18).Value = 7984.717648477152
$ws.Cells.Item(2, 19).Value = 0.04647783732217815
$ws.Cells.Item(2, 20).Value = 0.04647783732217814

# Row 3
$ws.Cells.Item(3, 7).Value = 51.08364133333333
$ws.Cells.Item(3, 8).Value = 153.250924
$ws.Cells.Item(3, 9).Value = 0.2324760785757066
$ws.Cells.Item(3, 10).Value = 0.2324760785757066
$ws.Cells.Item(3, 14).Value = 63.493212
$ws.Cells.Item(3, 15).Value = 0.2436343257635653
$ws.Cells.Item(3, 16).Value = 0.2436343257635653
$ws.Cells.Item(3, 17).Value = 1081.154822969765
$ws.Cells.Item(3, 18).Value = 9730.393406727888
$ws.Cells.Item(3, 19).Value = 0.05663915265994991
$ws.Cells.Item(3, 20).Value = 0.0566391526599499

# Row 4
$ws.Cells.Item(4, 7).Value = 51.08364133333333
$ws.Cells.Item(4, 8).Value = 153.250924
$ws.Cells.Item(4, 9).Value = 0.2324760785757066
$ws.Cells.Item(4, 10).Value = 0.2324760785757066
$ws.Cells.Item(4, 13).Value = 19.49944633333333
$ws.Cells.Item(4, 14).Value = 58.498339
$ws.Cells.Item(4, 15).Value = 0.2244681428394815
$ws.Cells.Item(4, 16).Value = 0.2244681428394814
$ws.Cells.Item(4, 17).Value = 996.1027226905818
$ws.Cells.Item(4, 18).Value = 8964.924504215236
$ws.Cells.Item(4, 19).Value = 0.05218347361249422
$ws.Cells.Item(4, 20).Value = 0.05218347361249421

# Row 5
$ws.Cells.Item(5, 7).Value = 51.08364133333333
$ws.Cells.Item(5, 8).Value = 153.250924
$ws.Cells.Item(5, 9).Value = 0.2324760785757066
$ws.Cells.Item(5, 10).Value = 0.2324760785757066
$ws.Cells.Item(5, 13).Value = 10.02612833333333
$ws.Cells.Item(5, 14).Value = 30.078385
$ws.Cells.Item(5, 15).Value = 0.1154159132716728
$ws.Cells.Item(5, 16).Value = 0.1154159132716728
$ws.Cells.Item(5, 17).Value = 512.171143741971
$ws.Cells.Item(5, 18).Value = 4609.540293677739
$ws.Cells.Item(5, 19).Value = 0.02683143892263235
$ws.Cells.Item(5, 20).Value = 0.02683143892263234

# Row 6
$ws.Cells.Item(6, 7).Value = 51.08364133333333
$ws.Cells.Item(6, 8).Value = 153.250924
$ws.Cells.Item(6, 9).Value = 0.2324760785757066
$ws.Cells.Item(6, 10).Value = 0.2324760785757066
$ws.Cells.Item(6, 13).Value = 18.81215433333334
$ws.Cells.Item(6, 14).Value = 56.436463
$ws.Cells.Item(6, 15).Value = 0.2165563715926893
$ws.Cells.Item(6, 16).Value = 0.2165563715926893
$ws.Cells.Item(6, 17).Value = 960.9933446713126
$ws.Cells.Item(6, 18).Value = 8648.940102041812
$ws.Cells.Item(6, 19).Value = 0.05034417605845196
$ws.Cells.Item(6, 20).Value = 0.05034417605845195

# Row 7
$ws.Cells.Item(7, 9).Value = 0.2058939713656488
$ws.Cells.Item(7, 10).Value = 0.2058939713656488
$ws.Cells.Item(7, 13).Value = 17.367416
$ws.Cells.Item(7, 14).Value = 52.102248
$ws.Cells.Item(7, 15).Value = 0.199925246532591
$ws.Cells.Item(7, 16).Value = 0.199925246532591
$ws.Cells.Item(7, 17).Value = 785.7464241054852
$ws.Cells.Item(7, 18).Value = 7071.717816949367
$ws.Cells.Item(7, 19).Value = 0.04116340298485158
$ws.Cells.Item(7, 20).Value = 0.04116340298485158

# Row 8
$ws.Cells.Item(8, 9).Value = 0.2058939713656488
$ws.Cells.Item(8, 10).Value = 0.2058939713656488
$ws.Cells.Item(8, 14).Value = 63.493212
$ws.Cells.Item(8, 15).Value = 0.2436343257635653
$ws.Cells.Item(8, 16).Value = 0.2436343257635653
$ws.Cells.Item(8, 18).Value = 8617.78705893349
$ws.Cells.Item(8, 19).Value = 0.05016283889245267
$ws.Cells.Item(8, 20).Value = 0.05016283889245267

# Row 9
$ws.Cells.Item(9, 9).Value = 0.2058939713656488
$ws.Cells.Item(9, 10).Value = 0.2058939713656488
$ws.Cells.Item(9, 13).Value = 19.49944633333333
$ws.Cells.Item(9, 14).Value = 58.498339
$ws.Cells.Item(9, 15).Value = 0.2244681428394815
$ws.Cells.Item(9, 16).Value = 0.2244681428394814
$ws.Cells.Item(9, 17).Value = 882.2049422005832
$ws.Cells.Item(9, 18).Value = 7939.844479805249
$ws.Cells.Item(9, 19).Value = 0.04621663737429256
$ws.Cells.Item(9, 20).Value = 0.04621663737429255

# Row 10
$ws.Cells.Item(10, 9).Value = 0.2058939713656488
$ws.Cells.Item(10, 10).Value = 0.2058939713656488
$ws.Cells.Item(10, 13).Value = 10.02612833333333
$ws.Cells.Item(10, 14).Value = 30.078385
$ws.Cells.Item(10, 15).Value = 0.1154159132716728
$ws.Cells.Item(10, 16).Value = 0.1154159132716728
$ws.Cells.Item(10, 17).Value = 453.6077494510038
$ws.Cells.Item(10, 18).Value = 4082.469745059034
$ws.Cells.Item(10, 19).Value = 0.02376344074229801
$ws.Cells.Item(10, 20).Value = 0.02376344074229801

# Row 11
$ws.Cells.Item(11, 9).Value = 0.2058939713656488
$ws.Cells.Item(11, 10).Value = 0.2058939713656488
$ws.Cells.Item(11, 13).Value = 18.81215433333334
$ws.Cells.Item(11, 14).Value = 56.436463
$ws.Cells.Item(11, 15).Value = 0.2165563715926893
$ws.Cells.Item(11, 16).Value = 0.2165563715926893
$ws.Cells.Item(11, 17).Value = 851.1100901329926
$ws.Cells.Item(11, 18).Value = 7659.990811196933
$ws.Cells.Item(11, 19).Value = 0.04458765137175398
$ws.Cells.Item(11, 20).Value = 0.04458765137175397

# Row 12
$ws.Cells.Item(12, 7).Value = 55.60882566666667
$ws.Cells.Item(12, 8).Value = 166.826477
$ws.Cells.Item(12, 9).Value = 0.253069698787332
$ws.Cells.Item(12, 10).Value = 0.253069698787332
$ws.Cells.Item(12, 13).Value = 17.367416
$ws.Cells.Item(12, 14).Value = 52.102248
$ws.Cells.Item(12, 15).Value = 0.199925246532591
$ws.Cells.Item(12, 16).Value = 0.199925246532591
$ws.Cells.Item(12, 17).Value = 965.7816086244773
$ws.Cells.Item(12, 18).Value = 8692.034477620296
$ws.Cells.Item(12, 19).Value = 0.0505950219199859
$ws.Cells.Item(12, 20).Value = 0.0505950219199859

# Row 13
$ws.Cells.Item(13, 7).Value = 55.60882566666667
$ws.Cells.Item(13, 8).Value = 166.826477
$ws.Cells.Item(13, 9).Value = 0.253069698787332
$ws.Cells.Item(13, 10).Value = 0.253069698787332
$ws.Cells.Item(13, 14).Value = 63.493212
$ws.Cells.Item(13, 15).Value = 0.2436343257635653
$ws.Cells.Item(13, 16).Value = 0.2436343257635653
$ws.Cells.Item(13, 17).Value = 1176.927652374903
$ws.Cells.Item(13, 18).Value = 10592.34887137412
$ws.Cells.Item(13, 19).Value = 0.06165646543524019
$ws.Cells.Item(13, 20).Value = 0.06165646543524019

# Row 14
$ws.Cells.Item(14, 7).Value = 55.60882566666667
$ws.Cells.Item(14, 8).Value = 166.826477
$ws.Cells.Item(14, 9).Value = 0.253069698787332
$ws.Cells.Item(14, 10).Value = 0.253069698787332
$ws.Cells.Item(14, 13).Value = 19.49944633333333
$ws.Cells.Item(14, 14).Value = 58.498339
$ws.Cells.Item(14, 15).Value = 0.2244681428394815
$ws.Cells.Item(14, 16).Value = 0.2244681428394814
$ws.Cells.Item(14, 17).Value = 1084.341311746856
$ws.Cells.Item(14, 18).Value = 9759.071805721704
$ws.Cells.Item(14, 19).Value = 0.05680608529573938
$ws.Cells.Item(14, 20).Value = 0.05680608529573938

# Row 15
$ws.Cells.Item(15, 7).Value = 55.60882566666667
$ws.Cells.Item(15, 8).Value = 166.826477
$ws.Cells.Item(15, 9).Value = 0.253069698787332
$ws.Cells.Item(15, 10).Value = 0.253069698787332
$ws.Cells.Item(15, 13).Value = 10.02612833333333
$ws.Cells.Item(15, 14).Value = 30.078385
$ws.Cells.Item(15, 15).Value = 0.1154159132716728
$ws.Cells.Item(15, 16).Value = 0.1154159132716728
$ws.Cells.Item(15, 17).Value = 557.5412225999605
$ws.Cells.Item(15, 18).Value = 5017.871003399645
$ws.Cells.Item(15, 19).Value = 0.02920827040692707
$ws.Cells.Item(15, 20).Value = 0.02920827040692707

# Row 16
$ws.Cells.Item(16, 7).Value = 55.60882566666667
$ws.Cells.Item(16, 8).Value = 166.826477
$ws.Cells.Item(16, 9).Value = 0.253069698787332
$ws.Cells.Item(16, 10).Value = 0.253069698787332
$ws.Cells.Item(16, 13).Value = 18.81215433333334
$ws.Cells.Item(16, 14).Value = 56.436463
$ws.Cells.Item(16, 15).Value = 0.2165563715926893
$ws.Cells.Item(16, 16).Value = 0.2165563715926893
$ws.Cells.Item(16, 17).Value = 1046.121810736761
$ws.Cells.Item(16, 18).Value = 9415.096296630853
$ws.Cells.Item(16, 19).Value = 0.05480385572943943
$ws.Cells.Item(16, 20).Value = 0.05480385572943942

# Row 17
$ws.Cells.Item(17, 7).Value = 25.04144866666667
$ws.Cells.Item(17, 8).Value = 75.124346
$ws.Cells.Item(17, 9).Value = 0.1139609009055278
$ws.Cells.Item(17, 10).Value = 0.1139609009055278
$ws.Cells.Item(17, 13).Value = 17.367416
$ws.Cells.Item(17, 14).Value = 52.102248
$ws.Cells.Item(17, 15).Value = 0.199925246532591
$ws.Cells.Item(17, 16).Value = 0.199925246532591
$ws.Cells.Item(17, 17).Value = 434.9052562366453
$ws.Cells.Item(17, 18).Value = 3914.147306129808
$ws.Cells.Item(17, 19).Value = 0.02278366120861382
$ws.Cells.Item(17, 20).Value = 0.02278366120861382

# Row 18
$ws.Cells.Item(18, 7).Value = 25.04144866666667
$ws.Cells.Item(18, 8).Value = 75.124346
$ws.Cells.Item(18, 9).Value = 0.1139609009055278
$ws.Cells.Item(18, 10).Value = 0.1139609009055278
$ws.Cells.Item(18, 14).Value = 63.493212
$ws.Cells.Item(18, 15).Value = 0.2436343257635653
$ws.Cells.Item(18, 16).Value = 0.2436343257635653
$ws.Cells.Item(18, 17).Value = 529.9873363265947
$ws.Cells.Item(18, 18).Value = 4769.886026939352
$ws.Cells.Item(18, 19).Value = 0.02776478725552674
$ws.Cells.Item(18, 20).Value = 0.02776478725552674

# Row 19
$ws.Cells.Item(19, 7).Value = 25.04144866666667
$ws.Cells.Item(19, 8).Value = 75.124346
$ws.Cells.Item(19, 9).Value = 0.1139609009055278
$ws.Cells.Item(19, 10).Value = 0.1139609009055278
$ws.Cells.Item(19, 13).Value = 19.49944633333333
$ws.Cells.Item(19, 14).Value = 58.498339
$ws.Cells.Item(19, 15).Value = 0.2244681428394815
$ws.Cells.Item(19, 16).Value = 0.2244681428394814
$ws.Cells.Item(19, 17).Value = 488.2943843845883
$ws.Cells.Item(19, 18).Value = 4394.649459461294
$ws.Cells.Item(19, 19).Value = 0.025580591782578
$ws.Cells.Item(19, 20).Value = 0.025580591782578

# Row 20
$ws.Cells.Item(20, 7).Value = 25.04144866666667
$ws.Cells.Item(20, 8).Value = 75.124346
$ws.Cells.Item(20, 9).Value = 0.1139609009055278
$ws.Cells.Item(20, 10).Value = 0.1139609009055278
$ws.Cells.Item(20, 13).Value = 10.02612833333333
$ws.Cells.Item(20, 14).Value = 30.078385
$ws.Cells.Item(20, 15).Value = 0.1154159132716728
$ws.Cells.Item(20, 16).Value = 0.1154159132716728
$ws.Cells.Item(20, 17).Value = 251.0687779845789
$ws.Cells.Item(20, 18).Value = 2259.61900186121
$ws.Cells.Item(20, 19).Value = 0.01315290145527409
$ws.Cells.Item(20, 20).Value = 0.01315290145527409

# Row 21
$ws.Cells.Item(21, 7).Value = 25.04144866666667
$ws.Cells.Item(21, 8).Value = 75.124346
$ws.Cells.Item(21, 9).Value = 0.1139609009055278
$ws.Cells.Item(21, 10).Value = 0.1139609009055278
$ws.Cells.Item(21, 13).Value = 18.81215433333334
$ws.Cells.Item(21, 14).Value = 56.436463
$ws.Cells.Item(21, 15).Value = 0.2165563715926893
$ws.Cells.Item(21, 16).Value = 0.2165563715926893
$ws.Cells.Item(21, 17).Value = 471.0835970475777
$ws.Cells.Item(21, 18).Value = 4239.752373428199
$ws.Cells.Item(21, 19).Value = 0.02467895920353512
$ws.Cells.Item(21, 20).Value = 0.02467895920353512

# Row 22
$ws.Cells.Item(22, 7).Value = 42.76071533333334
$ws.Cells.Item(22, 8).Value = 128.282146
$ws.Cells.Item(22, 9).Value = 0.1945993503657849
$ws.Cells.Item(22, 10).Value = 0.1945993503657849
$ws.Cells.Item(22, 13).Value = 17.367416
$ws.Cells.Item(22, 14).Value = 52.102248
$ws.Cells.Item(22, 15).Value = 0.199925246532591
$ws.Cells.Item(22, 16).Value = 0.199925246532591
$ws.Cells.Item(22, 17).Value = 742.6431316515786
$ws.Cells.Item(22, 18).Value = 6683.788184864208
$ws.Cells.Item(22, 19).Value = 0.0389053230969616
$ws.Cells.Item(22, 20).Value = 0.0389053230969616

# Row 23
$ws.Cells.Item(23, 7).Value = 42.76071533333334
$ws.Cells.Item(23, 8).Value = 128.282146
$ws.Cells.Item(23, 9).Value = 0.1945993503657849
$ws.Cells.Item(23, 10).Value = 0.1945993503657849
$ws.Cells.Item(23, 14).Value = 63.493212
$ws.Cells.Item(23, 15).Value = 0.2436343257635653
$ws.Cells.Item(23, 16).Value = 0.2436343257635653
$ws.Cells.Item(23, 17).Value = 905.0050546436614
$ws.Cells.Item(23, 18).Value = 8145.045491792953
$ws.Cells.Item(23, 19).Value = 0.04741108152039581
$ws.Cells.Item(23, 20).Value = 0.04741108152039581

# Row 24
$ws.Cells.Item(24, 7).Value = 42.76071533333334
$ws.Cells.Item(24, 8).Value = 128.282146
$ws.Cells.Item(24, 9).Value = 0.1945993503657849
$ws.Cells.Item(24, 10).Value = 0.1945993503657849
$ws.Cells.Item(24, 13).Value = 19.49944633333333
$ws.Cells.Item(24, 14).Value = 58.498339
$ws.Cells.Item(24, 15).Value = 0.2244681428394815
$ws.Cells.Item(24, 16).Value = 0.2244681428394814
$ws.Cells.Item(24, 17).Value = 833.8102738172772
$ws.Cells.Item(24, 18).Value = 7504.292464355495
$ws.Cells.Item(24, 19).Value = 0.0436813547743773
$ws.Cells.Item(24, 20).Value = 0.04368135477437729

# Row 25
$ws.Cells.Item(25, 7).Value = 42.76071533333334
$ws.Cells.Item(25, 8).Value = 128.282146
$ws.Cells.Item(25, 9).Value = 0.1945993503657849
$ws.Cells.Item(25, 10).Value = 0.1945993503657849
$ws.Cells.Item(25, 13).Value = 10.02612833333333
$ws.Cells.Item(25, 14).Value = 30.078385
$ws.Cells.Item(25, 15).Value = 0.1154159132716728
$ws.Cells.Item(25, 16).Value = 0.1154159132716728
$ws.Cells.Item(25, 17).Value = 428.7244195571345
$ws.Cells.Item(25, 18).Value = 3858.51977601421
$ws.Cells.Item(25, 19).Value = 0.0224598617445413
$ws.Cells.Item(25, 20).Value = 0.0224598617445413

# Row 26
$ws.Cells.Item(26, 7).Value = 42.76071533333334
$ws.Cells.Item(26, 8).Value = 128.282146
$ws.Cells.Item(26, 9).Value = 0.1945993503657849
$ws.Cells.Item(26, 10).Value = 0.1945993503657849
$ws.Cells.Item(26, 13).Value = 18.81215433333334
$ws.Cells.Item(26, 14).Value = 56.436463
$ws.Cells.Item(26, 15).Value = 0.2165563715926893
$ws.Cells.Item(26, 16).Value = 0.2165563715926893
$ws.Cells.Item(26, 17).Value = 804.4211762543999
$ws.Cells.Item(26, 18).Value = 7239.790586289599
$ws.Cells.Item(26, 19).Value = 0.04214172922950885
$ws.Cells.Item(26, 20).Value = 0.04214172922950885
